# Update project plan worksheet:
#  - Clear the "请假" (leave of absence) note for 苏立明 on the first week's table.
#  - Insert a new row for 蔡智杰 in the second week's table (before the "总结：" block)
#    and fill in the member/plan columns for that whole second table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Clear the remark for 苏立明 (row 8, column D) in the first table.
$ws.Range("D8").Value = ""

# 2) Insert a new row above the old row 18 (the merged "总结：" block), which
#    shifts that merged block down from A18:D19 to A19:D20 and makes room for
#    a sixth team member's row in the second table.
$ws.Cells.Item(18, 1).EntireRow.Insert()

# Copy the formatting of the row above (row 17, a normal data row) onto the
# freshly inserted row 18 so it matches the rest of the table's borders/font.
$ws.Range("A17:D17").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Fill in the member names and plan content for the second table
#    (rows 13-18), now that row 18 exists for 蔡智杰.
$ws.Range("A13").Value = "何舒静"
$ws.Range("B13").Value = "商讨，修改完善用例图，合并用例图"

$ws.Range("A14").Value = "陈碧容"
$ws.Range("B14").Value = "商讨，修改完善用例图"

$ws.Range("A15").Value = "黄丙升"
$ws.Range("B15").Value = "商讨，修改完善用例图"

$ws.Range("A16").Value = "王增璟"
$ws.Range("B16").Value = "商讨，修改完善用例图"

$ws.Range("A17").Value = "苏立明"
$ws.Range("B17").Value = "商讨，修改完善用例图"

$ws.Range("A18").Value = "蔡智杰"
$ws.Range("B18").Value = "商讨，修改完善用例图"

# 4) Match the author's final selection/active cell.
$ws.Range("D18").Select()
